$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(8, 9).Value = 'aa'
$ws.Cells.Item(8, 10).Value = 'Agree/Accept'
$ws.Cells.Item(11, 9).Value = 'sd'
$ws.Cells.Item(11, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(15, 9).Value = 'b'
$ws.Cells.Item(15, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(16, 9).Value = 'sd'
$ws.Cells.Item(16, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(17, 9).Value = 'aa'
$ws.Cells.Item(17, 10).Value = 'Agree/Accept'
$ws.Cells.Item(23, 9).Value = 'b'
$ws.Cells.Item(23, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(27, 9).Value = 'aa'
$ws.Cells.Item(27, 10).Value = 'Agree/Accept'
$ws.Cells.Item(49, 9).Value = 'ba'
$ws.Cells.Item(49, 10).Value = 'Appreciation'
$ws.Cells.Item(50, 9).Value = '%'
$ws.Cells.Item(50, 10).Value = 'Uninterpretable'
$ws.Cells.Item(64, 9).Value = 'b'
$ws.Cells.Item(64, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(65, 9).Value = 'qy'
$ws.Cells.Item(65, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(71, 9).Value = 'b'
$ws.Cells.Item(71, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(86, 9).Value = 'aa'
$ws.Cells.Item(86, 10).Value = 'Agree/Accept'
$ws.Cells.Item(93, 9).Value = 'b'
$ws.Cells.Item(93, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(94, 9).Value = 'aa'
$ws.Cells.Item(94, 10).Value = 'Agree/Accept'
$ws.Cells.Item(110, 9).Value = 'sv'
$ws.Cells.Item(110, 10).Value = 'Statement-opinion'
$ws.Cells.Item(114, 9).Value = 'sd'
$ws.Cells.Item(114, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(124, 9).Value = 'sv'
$ws.Cells.Item(124, 10).Value = 'Statement-opinion'
$ws.Cells.Item(142, 9).Value = 'sv'
$ws.Cells.Item(142, 10).Value = 'Statement-opinion'
$ws.Cells.Item(143, 9).Value = 'sv'
$ws.Cells.Item(143, 10).Value = 'Statement-opinion'
$ws.Cells.Item(146, 9).Value = 'sv'
$ws.Cells.Item(146, 10).Value = 'Statement-opinion'
$ws.Cells.Item(151, 9).Value = 'aa'
$ws.Cells.Item(151, 10).Value = 'Agree/Accept'
$ws.Cells.Item(158, 9).Value = 'aa'
$ws.Cells.Item(158, 10).Value = 'Agree/Accept'
$ws.Cells.Item(160, 9).Value = 'sd'
$ws.Cells.Item(160, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(164, 9).Value = 'sv'
$ws.Cells.Item(164, 10).Value = 'Statement-opinion'
$ws.Cells.Item(166, 9).Value = 'sd'
$ws.Cells.Item(166, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(180, 9).Value = 'ba'
$ws.Cells.Item(180, 10).Value = 'Appreciation'
$ws.Cells.Item(181, 9).Value = 'qy'
$ws.Cells.Item(181, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(182, 9).Value = 'sv'
$ws.Cells.Item(182, 10).Value = 'Statement-opinion'
$ws.Cells.Item(204, 9).Value = 'sv'
$ws.Cells.Item(204, 10).Value = 'Statement-opinion'
$ws.Cells.Item(206, 9).Value = 'sd'
$ws.Cells.Item(206, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(217, 9).Value = 'aa'
$ws.Cells.Item(217, 10).Value = 'Agree/Accept'
$ws.Cells.Item(218, 9).Value = 'sv'
$ws.Cells.Item(218, 10).Value = 'Statement-opinion'
$ws.Cells.Item(219, 9).Value = 'sv'
$ws.Cells.Item(219, 10).Value = 'Statement-opinion'
$ws.Cells.Item(220, 9).Value = 'aa'
$ws.Cells.Item(220, 10).Value = 'Agree/Accept'
$ws.Cells.Item(223, 9).Value = 'b'
$ws.Cells.Item(223, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(235, 9).Value = 'sv'
$ws.Cells.Item(235, 10).Value = 'Statement-opinion'
$ws.Cells.Item(236, 9).Value = 'sd'
$ws.Cells.Item(236, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(238, 9).Value = 'sv'
$ws.Cells.Item(238, 10).Value = 'Statement-opinion'
$ws.Cells.Item(240, 9).Value = 'sd'
$ws.Cells.Item(240, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(241, 9).Value = 'aa'
$ws.Cells.Item(241, 10).Value = 'Agree/Accept'
$ws.Cells.Item(247, 9).Value = 'aa'
$ws.Cells.Item(247, 10).Value = 'Agree/Accept'
$ws.Cells.Item(249, 9).Value = 'b'
$ws.Cells.Item(249, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(252, 9).Value = 'sv'
$ws.Cells.Item(252, 10).Value = 'Statement-opinion'
$ws.Cells.Item(256, 9).Value = 'sv'
$ws.Cells.Item(256, 10).Value = 'Statement-opinion'
$ws.Cells.Item(267, 9).Value = 'aa'
$ws.Cells.Item(267, 10).Value = 'Agree/Accept'
$ws.Cells.Item(268, 9).Value = 'sv'
$ws.Cells.Item(268, 10).Value = 'Statement-opinion'
$ws.Cells.Item(271, 9).Value = 'b'
$ws.Cells.Item(271, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(274, 9).Value = 'aa'
$ws.Cells.Item(274, 10).Value = 'Agree/Accept'
$ws.Cells.Item(287, 9).Value = 'aa'
$ws.Cells.Item(287, 10).Value = 'Agree/Accept'
$ws.Cells.Item(288, 9).Value = 'ba'
$ws.Cells.Item(288, 10).Value = 'Appreciation'
$ws.Cells.Item(291, 9).Value = '%'
$ws.Cells.Item(291, 10).Value = 'Uninterpretable'
$ws.Cells.Item(293, 9).Value = 'aa'
$ws.Cells.Item(293, 10).Value = 'Agree/Accept'
$ws.Cells.Item(309, 9).Value = 'qy'
$ws.Cells.Item(309, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(315, 9).Value = 'sv'
$ws.Cells.Item(315, 10).Value = 'Statement-opinion'
$ws.Cells.Item(324, 9).Value = 'sd'
$ws.Cells.Item(324, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(325, 9).Value = 'sd'
$ws.Cells.Item(325, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(327, 9).Value = 'sd'
$ws.Cells.Item(327, 10).Value = 'Statement-non-opinion'
